$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps storing values as text (matches the
# original workbook, where prices like "30.344.74" are plain strings,
# not numbers) so Excel does not auto-convert numeric-looking values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.351.46"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.870.15"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "235.62"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.2850"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +6.27%  "
$ws.Range("D11").Value = "0.07884"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "98.32"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "1.873.44"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "0.6768"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "278.65"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "30.350.39"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "5.483"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "2.120.95"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "0.000007309"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "6.159"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "165.48"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").Value = "9.169"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").Value = "19.19"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").Value = "1.935"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "1.381"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").Value = "0.09700"
$ws.Range("D31").Value = "4.401"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").Value = "0.04714"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D35").Value = "1.130"
$ws.Range("E35").Value = "  +4.19%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").Value = "6.321"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "2.541"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "74.70"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").Value = "1.954"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "0.8506"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "0.4191"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "103.98"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "7.228"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").Value = "937.77"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("D50").Value = "34.25"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "0.1123"
$ws.Range("E51").Value = "  -1.50%  "
